# Add "Increase Tide" and "Decrease Tide" reward/powerup rows to the
# localization sheet. This inserts two new rows right before the existing
# "LESSON_47" row (old row 32, now row 34) so the new powerup strings sit
# with the other UI_REWARDS_* entries, and shifts everything below down by
# two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 32 (pushes old row 32.. down to 34..).
$ws.Rows.Item(32).EntireRow.Insert()
$ws.Rows.Item(32).EntireRow.Insert()

# Row 31 (A31:E31) already carries the exact style combination the new
# rows need (s=7,12,7,8,12 with the 15pt "thickBot" row look), so copy its
# formatting down into the two freshly inserted rows.
$ws.Range("A31:E31").Copy()
$ws.Range("A32:E33").PasteSpecial(-4122)
$ws.Rows.Item(32).RowHeight = 15
$ws.Rows.Item(33).RowHeight = 15

# Fill in the Key column for both new rows first (so the two new keys land
# in the shared-string table before the two new English values do).
$ws.Cells.Item(32, 1).Value = "UI_REWARDS_DECREASE_TIDE"
$ws.Cells.Item(33, 1).Value = "UI_REWARDS_INCREASE_TIDE"

# Row 33 (Increase Tide): English text, then the still-untranslated columns.
$ws.Cells.Item(33, 2).Value = "Water Speed Increased"
$ws.Cells.Item(33, 3).Value = "XXXX"
$ws.Cells.Item(33, 4).Value = "XXXX"
$ws.Cells.Item(33, 5).Value = "XXXX"

# Row 32 (Decrease Tide): English text, then the still-untranslated columns.
$ws.Cells.Item(32, 2).Value = "Water Speed Decreased"
$ws.Cells.Item(32, 3).Value = "XXXX"
$ws.Cells.Item(32, 4).Value = "XXXX"
$ws.Cells.Item(32, 5).Value = "XXXX"

# Restore the view to where it was when the edit was made.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B33").Select()
